$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 13.62572784085251
$ws.Range("D2").Value = 8.862755317428544
$ws.Range("E2").Value = 14.33113577209349
$ws.Range("F2").Value = 37.46470187668112
$ws.Range("G2").Value = 42.83247596206252
$ws.Range("H2").Value = 17.26068741337965
$ws.Range("J2").Value = 10.70174230392237
$ws.Range("L2").Value = 9.757431744702821
$ws.Range("N2").Value = 18.99769123883691
$ws.Range("O2").Value = 28.31248297128701

$ws.Range("C3").Value = 13.59124863422203
$ws.Range("D3").Value = 8.873921610170255
$ws.Range("E3").Value = 14.34229924205262
$ws.Range("F3").Value = 37.26013579872049
$ws.Range("G3").Value = 42.33423386470002
$ws.Range("H3").Value = 17.23884865881928
$ws.Range("J3").Value = 10.71951035084092
$ws.Range("L3").Value = 9.775400896148135
$ws.Range("N3").Value = 18.40031508502705
$ws.Range("O3").Value = 28.18243555185087

$ws.Range("C4").Value = 13.57300577819324
$ws.Range("D4").Value = 8.881897797228048
$ws.Range("E4").Value = 14.35144590203818
$ws.Range("F4").Value = 37.14330426916217
$ws.Range("G4").Value = 42.03776962370898
$ws.Range("H4").Value = 17.22888948356952
$ws.Range("J4").Value = 10.73176849844094
$ws.Range("L4").Value = 9.787315083032775
$ws.Range("N4").Value = 18.02485520896364
$ws.Range("O4").Value = 28.10908694031781

$ws.Range("C5").Value = 13.56631267829622
$ws.Range("D5").Value = 8.885430033877224
$ws.Range("E5").Value = 14.35574987205538
$ws.Range("F5").Value = 37.09793849221111
$ws.Range("G5").Value = 41.91947647281248
$ws.Range("H5").Value = 17.2257013727123
$ws.Range("J5").Value = 10.73710307738097
$ws.Range("L5").Value = 9.792392179752509
$ws.Range("N5").Value = 17.86990355188765
$ws.Range("O5").Value = 28.08085364041378

$ws.Range("C6").Value = 13.56524617732769
$ws.Range("D6").Value = 8.88603358990474
$ws.Range("E6").Value = 14.35649937399993
$ws.Range("F6").Value = 37.09054211389126
$ws.Range("G6").Value = 41.89999007353464
$ws.Range("H6").Value = 17.22522462595809
$ws.Range("J6").Value = 10.73800937736525
$ws.Range("L6").Value = 9.793248646511472
$ws.Range("N6").Value = 17.84406337566597
$ws.Range("O6").Value = 28.076266192773

$ws.Range("C7").Value = 13.57291250648684
$ws.Range("D7").Value = 8.881944292612792
$ws.Range("E7").Value = 14.35150161192658
$ws.Range("F7").Value = 37.14268331700409
$ws.Range("G7").Value = 42.0361639031463
$ws.Range("H7").Value = 17.22884296042277
$ws.Range("J7").Value = 10.73183906845057
$ws.Range("L7").Value = 9.787382655234515
$ws.Range("N7").Value = 18.02277304767602
$ws.Range("O7").Value = 28.10869944022903

$ws.Range("C8").Value = 13.61323493729686
$ws.Range("D8").Value = 8.866373138898705
$ws.Range("E8").Value = 14.3345093573034
$ws.Range("F8").Value = 37.39237054819623
$ws.Range("G8").Value = 42.65880465962591
$ws.Range("H8").Value = 17.25244325448834
$ws.Range("J8").Value = 10.70758896840152
$ws.Range("L8").Value = 9.763444897773967
$ws.Range("N8").Value = 18.79364780656866
$ws.Range("O8").Value = 28.26630730001685

$ws.Range("C9").Value = 13.71529058171047
$ws.Range("D9").Value = 8.844714389750035
$ws.Range("E9").Value = 14.31936133740966
$ws.Range("F9").Value = 37.94978392962692
$ws.Range("G9").Value = 43.94845170283357
$ws.Range("H9").Value = 17.32594814660332
$ws.Range("J9").Value = 10.67072554505185
$ws.Range("L9").Value = 9.723475175503433
$ws.Range("N9").Value = 20.2273683202997
$ws.Range("O9").Value = 28.62588138305023

$ws.Range("C10").Value = 13.80390930595297
$ws.Range("D10").Value = 8.834197609021146
$ws.Range("E10").Value = 14.31928596505052
$ws.Range("F10").Value = 38.39796505636185
$ws.Range("G10").Value = 44.92898426115371
$ws.Range("H10").Value = 17.39631425799098
$ws.Range("J10").Value = 10.65014734746823
$ws.Range("L10").Value = 9.698334533188332
$ws.Range("N10").Value = 21.22223697909767
$ws.Range("O10").Value = 28.91926276991264

$ws.Range("C11").Value = 13.84709189678716
$ws.Range("D11").Value = 8.830581042242571
$ws.Range("E11").Value = 14.32164344846591
$ws.Range("F11").Value = 38.6096184240628
$ws.Range("G11").Value = 45.38031206323659
$ws.Range("H11").Value = 17.43181581195757
$ws.Range("J11").Value = 10.64219537716486
$ws.Range("L11").Value = 9.687809496925762
$ws.Range("N11").Value = 21.66018057919901
$ws.Range("O11").Value = 29.05869084702697

$ws.Range("C12").Value = 13.86384765265519
$ws.Range("D12").Value = 8.829379038811584
$ws.Range("E12").Value = 14.32287906416171
$ws.Range("F12").Value = 38.69082727827901
$ws.Range("G12").Value = 45.55181696055553
$ws.Range("H12").Value = 17.44575508117021
$ws.Range("J12").Value = 10.63938650960049
$ws.Range("L12").Value = 9.683954600181984
$ws.Range("N12").Value = 21.82377585682186
$ws.Range("O12").Value = 29.11231200333785

$ws.Range("C13").Value = 13.86022121028739
$ws.Range("D13").Value = 8.829630468129295
$ws.Range("E13").Value = 14.32259772063719
$ws.Range("F13").Value = 38.6732912775461
$ws.Range("G13").Value = 45.51485633245201
$ws.Range("H13").Value = 17.44273107786221
$ws.Range("J13").Value = 10.63998245347801
$ws.Range("L13").Value = 9.684779013964178
$ws.Range("N13").Value = 21.78864458690801
$ws.Range("O13").Value = 29.10072771395259

$ws.Range("C14").Value = 13.84846237205802
$ws.Range("D14").Value = 8.830478797838275
$ws.Range("E14").Value = 14.32173823792006
$ws.Range("F14").Value = 38.61627861159151
$ws.Range("G14").Value = 45.39441083596279
$ws.Range("H14").Value = 17.43295271075696
$ws.Range("J14").Value = 10.64196023573319
$ws.Range("L14").Value = 9.687489734858826
$ws.Range("N14").Value = 21.67368539489659
$ws.Range("O14").Value = 29.06308601017968

$ws.Range("C15").Value = 13.8413120015017
$ws.Range("D15").Value = 8.831020227254879
$ws.Range("E15").Value = 14.32125640013918
$ws.Range("F15").Value = 38.58149305073947
$ws.Range("G15").Value = 45.32070743187435
$ws.Range("H15").Value = 17.42702751753759
$ws.Range("J15").Value = 10.64319803084872
$ws.Range("L15").Value = 9.689167140311888
$ws.Range("N15").Value = 21.60297336126124
$ws.Range("O15").Value = 29.04013545344801

$ws.Range("C16").Value = 13.8011441929024
$ws.Range("D16").Value = 8.834457422765899
$ws.Range("E16").Value = 14.319179940967
$ws.Range("F16").Value = 38.38428456588802
$ws.Range("G16").Value = 44.8995812495741
$ws.Range("H16").Value = 17.39406386747884
$ws.Range("J16").Value = 10.65069536131443
$ws.Range("L16").Value = 9.699040679950887
$ws.Range("N16").Value = 21.19330956972086
$ws.Range("O16").Value = 28.91026791198713

$ws.Range("C17").Value = 13.77723115553604
$ws.Range("D17").Value = 8.836864812886526
$ws.Range("E17").Value = 14.31851799987657
$ws.Range("F17").Value = 38.26525429078887
$ws.Range("G17").Value = 44.64247014623312
$ws.Range("H17").Value = 17.37473160250478
$ws.Range("J17").Value = 10.65565547018049
$ws.Range("L17").Value = 9.705330983097408
$ws.Range("N17").Value = 20.93814219015166
$ws.Range("O17").Value = 28.83210214360583

$ws.Range("C18").Value = 13.76374768952227
$ws.Range("D18").Value = 8.838359402522137
$ws.Range("E18").Value = 14.31836246696413
$ws.Range("F18").Value = 38.19752649068086
$ws.Range("G18").Value = 44.49509314438173
$ws.Range("H18").Value = 17.36394127360019
$ws.Range("J18").Value = 10.65864104749173
$ws.Range("L18").Value = 9.709034822654257
$ws.Range("N18").Value = 20.79000725568362
$ws.Range("O18").Value = 28.78770716343064

$ws.Range("C19").Value = 13.75922917870536
$ws.Range("D19").Value = 8.83888433523278
$ws.Range("E19").Value = 14.31834851326144
$ws.Range("F19").Value = 38.17472302161086
$ws.Range("G19").Value = 44.44528578053139
$ws.Range("H19").Value = 17.36034456122519
$ws.Range("J19").Value = 10.65967470367349
$ws.Range("L19").Value = 9.710303630839036
$ws.Range("N19").Value = 20.73962067985785
$ws.Range("O19").Value = 28.77277368737888

$ws.Range("C20").Value = 13.77974878931552
$ws.Range("D20").Value = 8.836597168198232
$ws.Range("E20").Value = 14.31856516455195
$ws.Range("F20").Value = 38.27784957978086
$ws.Range("G20").Value = 44.66978882063518
$ws.Range("H20").Value = 17.37675554200244
$ws.Range("J20").Value = 10.65511373079685
$ws.Range("L20").Value = 9.704652489920084
$ws.Range("N20").Value = 20.96544799484619
$ws.Range("O20").Value = 28.84036491251916

$ws.Range("C21").Value = 13.85190535474531
$ws.Range("D21").Value = 8.830225079856596
$ws.Range("E21").Value = 14.32198139212292
$ws.Range("F21").Value = 38.63299630516158
$ws.Range("G21").Value = 45.429773684286
$ws.Range("H21").Value = 17.43581145843775
$ws.Range("J21").Value = 10.64137382309127
$ws.Range("L21").Value = 9.686689986092576
$ws.Range("N21").Value = 21.70751365554066
$ws.Range("O21").Value = 29.07412025053724

$ws.Range("C22").Value = 13.90141074261689
$ws.Range("D22").Value = 8.827036734604075
$ws.Range("E22").Value = 14.3262121026306
$ws.Range("F22").Value = 38.87125753551921
$ws.Range("G22").Value = 45.92987274775883
$ws.Range("H22").Value = 17.4772932284888
$ws.Range("J22").Value = 10.63357343013874
$ws.Range("L22").Value = 9.675712132454317
$ws.Range("N22").Value = 22.17935961385674
$ws.Range("O22").Value = 29.23167117122071

$ws.Range("C23").Value = 13.87477723634465
$ws.Range("D23").Value = 8.828649230142597
$ws.Range("E23").Value = 14.32377166416428
$ws.Range("F23").Value = 38.74354959853643
$ws.Range("G23").Value = 45.66270203397093
$ws.Range("H23").Value = 17.45489191910592
$ws.Range("J23").Value = 10.63762882299493
$ws.Range("L23").Value = 9.681501649785304
$ws.Range("N23").Value = 21.92877110912573
$ws.Range("O23").Value = 29.14715812794686

$ws.Range("C24").Value = 13.77860974266206
$ws.Range("D24").Value = 8.836717825994862
$ws.Range("E24").Value = 14.31854314040878
$ws.Range("F24").Value = 38.27215305146716
$ws.Range("G24").Value = 44.65743666625374
$ws.Range("H24").Value = 17.37583950907096
$ws.Range("J24").Value = 10.65535823393436
$ws.Range("L24").Value = 9.704958964230631
$ws.Range("N24").Value = 20.95310750188672
$ws.Range("O24").Value = 28.83662761886639

$ws.Range("C25").Value = 13.68525697606356
$ws.Range("D25").Value = 8.849624834287869
$ws.Range("E25").Value = 14.32151563862589
$ws.Range("F25").Value = 37.79199649423815
$ws.Range("G25").Value = 43.59311989568475
$ws.Range("H25").Value = 17.30317030201022
$ws.Range("J25").Value = 10.67955469390841
$ws.Range("L25").Value = 9.733544246812025
$ws.Range("N25").Value = 19.84905939529495
$ws.Range("O25").Value = 28.52335573645843
